# Planning.xlsx update "on 11th of july"
# Fills in the newly-added weekly-planner entries and restyles the two
# cells that got a wrap + highlight treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---------------------------------------------------------------
$ws.Range("F5").Value = "Geo BBQ"
$ws.Range("G5").Value = "Annet Weekend"

# --- Row 6 -----------------------------------------------------------------
$ws.Range("B6").Value = "Weekly Review doen, RQT graphs van SSL Slam en ORB Slam maken, ROS doorspitten"
$ws.Range("B6").WrapText = $true

$ws.Range("C6").Value = "Onderzoeksvraag maken, Carlas en Robert Update sturen"
$ws.Range("C6").WrapText = $true
$ws.Range("C6").Interior.Color = 7561697

$ws.Range("F6").Value = "Progress Meeting"
$ws.Range("G6").Value = "TI Spelen"
$ws.Range("H6").Value = "Naar Ouders"

# --- Row 7 -------------------------------------------------------------
$ws.Range("G7").Value = "Pre Vakantie Weekend"

# --- Row 10 (Bruiloft text unchanged, kept for completeness) -------------
$ws.Range("F10").Value = "Bruiloft"

# Selection moved to C6 in the saved file.
$ws.Range("C6").Select()
